# Add a new "Save" column (H) to the sheet:
#  - H1 header "Save", formatted like the other header cells (reuse style of G1)
#  - H2:H4 numeric flag values (era data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell with the same look & feel as the existing header row.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for the "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
